$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows 6-9 entirely so the sheet shrinks to A1:B5
$ws.Range("A6:B9").EntireRow.Delete()

# Update the surviving rows (2-5) with their new values
$ws.Range("A2").Value = 21
$ws.Range("B2").Value = 38

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 36

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 22

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 3
